$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new sheets in the right places.
#    NOTE: sheet references returned/looked-up earlier become stale (they are
#    resolved by position, not identity) once the sheet collection is
#    mutated, so we always re-fetch by name right before using them.
# ---------------------------------------------------------------------------

# "Borrow" goes directly before "Electronics"
$borrow = $wb.Worksheets.Add($wb.Worksheets.Item("Electronics"))
$borrow.Name = "Borrow"

# "Users" goes directly before "Borrow" (so final order starts Users, Borrow, Electronics, ...)
$users = $wb.Worksheets.Add($wb.Worksheets.Item("Borrow"))
$users.Name = "Users"

# "Statestieken" goes at the very end, after "Books"
$stats = $wb.Worksheets.Add($null, $wb.Worksheets.Item("Books"))
$stats.Name = "Statestieken"

# ---------------------------------------------------------------------------
# 2. Populate "Users" sheet.
#    Values are written in the same order the shared strings were first
#    introduced so the rebuilt sharedStrings.xml lines up with the source.
# ---------------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A1").Value = "_ID"
$wsUsers.Range("B1").Value = "FIRSTNAME"
$wsUsers.Range("C1").Value = "SURNAME"
$wsUsers.Range("D1").Value = "SCHOOL_EMAIL"
$wsUsers.Range("E1").Value = "PASSWORD"
$wsUsers.Range("F1").Value = "USER_TYPE"
$wsUsers.Range("G1").Value = "BLOCKED"
$wsUsers.Range("G2").Value = "boolean"
$wsUsers.Range("F2").Value = "student"
$wsUsers.Range("F3").Value = "admin"
$wsUsers.Range("F4").Value = "beheer"
$wsUsers.Range("A2").Value = "INT"
$wsUsers.Range("H1").Value = "CLASSCODE"

# ---------------------------------------------------------------------------
# 3. Populate "Borrow" sheet.
# ---------------------------------------------------------------------------
$wsBorrow = $wb.Worksheets.Item("Borrow")
$wsBorrow.Range("A1").Value = "_ID"
$wsBorrow.Range("C1").Value = "P_USER_ID"
$wsBorrow.Range("B1").Value = "P_PRODUCT_ID"
$wsBorrow.Range("D1").Value = "DATETIME"
$wsBorrow.Range("A2").Value = "Int"
$wsBorrow.Range("F1").Value = "STATUS"
$wsBorrow.Range("E1").Value = "AMOUNT"
$wsBorrow.Range("B2").Value = "Int"
$wsBorrow.Range("C2").Value = "Int"
$wsBorrow.Range("D2").Value = "Int"
$wsBorrow.Range("E2").Value = "int"
$wsBorrow.Range("F2").Value = "text"

# ---------------------------------------------------------------------------
# 4. Column widths (best effort - engine quantises to 1/6-character pixels).
# ---------------------------------------------------------------------------
$wsUsers.Columns.Item(2).ColumnWidth = 9.3333333333
$wsUsers.Columns.Item(3).ColumnWidth = 8.3333333333
$wsUsers.Columns.Item(4).ColumnWidth = 13.1666666667
$wsUsers.Columns.Item(5).ColumnWidth = 13.1666666667
$wsUsers.Columns.Item(6).ColumnWidth = 9.3333333333
$wsUsers.Columns.Item(7).ColumnWidth = 9.3333333333
$wsUsers.Columns.Item(9).ColumnWidth = 8.3333333333

$wsBorrow.Columns.Item(2).ColumnWidth = 13
$wsBorrow.Columns.Item(3).ColumnWidth = 9

# ---------------------------------------------------------------------------
# 5. Selections / active cells to match the source workbook.
# ---------------------------------------------------------------------------
$wsUsers.Range("D14").Select()
$wsBorrow.Range("F2").Select()
$wb.Worksheets.Item("Electronics").Range("D9").Select()

$wsUsers.Select()

Write-Output "Sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
